# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Cactuar_Profits workbook (46 rows across 8 sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5,8).Value = 1385.375  # H5: 1540.4286 -> 1385.375
$ws.Cells.Item(5,9).Value = 181.5  # I5: 157.8 -> 181.5
$ws.Cells.Item(5,11).Value = 181.5  # K5: 157.8 -> 181.5
$ws.Cells.Item(5,13).Value = -66.5  # M5: -42.80000000000001 -> -66.5

$ws.Cells.Item(101,8).Value = 2762.4707  # H101: 2762.5293 -> 2762.4707
$ws.Cells.Item(101,9).Value = 1302.6364  # I101: 1302.7273 -> 1302.6364
$ws.Cells.Item(101,11).Value = 3907.9092  # K101: 3908.1819 -> 3907.9092
$ws.Cells.Item(101,13).Value = -2285.9092  # M101: -2286.1819 -> -2285.9092

$ws.Cells.Item(111,8).Value = 7050  # H111: 7335.7144 -> 7050
$ws.Cells.Item(111,9).Value = 9237.75  # I111: 11317.333 -> 9237.75
$ws.Cells.Item(111,10).Value = 4133  # J111: 4349.5 -> 4133
$ws.Cells.Item(111,11).Value = 27713.25  # K111: 33951.999 -> 27713.25
$ws.Cells.Item(111,12).Value = 12399  # L111: 13048.5 -> 12399
$ws.Cells.Item(111,13).Value = -24646.25  # M111: -30884.999 -> -24646.25
$ws.Cells.Item(111,14).Value = -18533  # N111: -19182.5 -> -18533

$ws.Cells.Item(132,8).Value = 316793  # H132: 417212.72 -> 316793
$ws.Cells.Item(132,9).Value = 427113.62  # I132: 640076.4 -> 427113.62
$ws.Cells.Item(132,11).Value = 1281340.86  # K132: 1920229.2 -> 1281340.86
$ws.Cells.Item(132,13).Value = -1278810.86  # M132: -1917699.2 -> -1278810.86

$ws.Cells.Item(137,8).Value = 2886.0789  # H137: 2868.487 -> 2886.0789
$ws.Cells.Item(137,9).Value = 1750.5  # I137: 1640.4 -> 1750.5
$ws.Cells.Item(137,10).Value = 3019.6765  # J137: 3049.0881 -> 3019.6765
$ws.Cells.Item(137,11).Value = 5251.5  # K137: 4921.200000000001 -> 5251.5
$ws.Cells.Item(137,12).Value = 9059.029500000001  # L137: 9147.264299999999 -> 9059.029500000001
$ws.Cells.Item(137,13).Value = -2701.5  # M137: -2371.200000000001 -> -2701.5
$ws.Cells.Item(137,14).Value = -14159.0295  # N137: -14247.2643 -> -14159.0295

$ws.Cells.Item(138,8).Value = 3544.42  # H138: 3714.94 -> 3544.42
$ws.Cells.Item(138,10).Value = 4088.6301  # J138: 4322.219 -> 4088.6301
$ws.Cells.Item(138,12).Value = 12265.8903  # L138: 12966.657 -> 12265.8903
$ws.Cells.Item(138,14).Value = -22545.8903  # N138: -23246.657 -> -22545.8903

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122,8).Value = 3851.1292  # H122: 4692.1665 -> 3851.1292
$ws.Cells.Item(122,9).Value = 2823.1  # I122: 3656.3572 -> 2823.1
$ws.Cells.Item(122,10).Value = 5720.273  # J122: 6142.3 -> 5720.273
$ws.Cells.Item(122,11).Value = 8469.299999999999  # K122: 10969.0716 -> 8469.299999999999
$ws.Cells.Item(122,12).Value = 17160.819  # L122: 18426.9 -> 17160.819
$ws.Cells.Item(122,13).Value = -6019.299999999999  # M122: -8519.071599999999 -> -6019.299999999999
$ws.Cells.Item(122,14).Value = -22060.819  # N122: -23326.9 -> -22060.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82,8).Value = 14336.429  # H82: 13044.375 -> 14336.429
$ws.Cells.Item(82,9).Value = 4725.8335  # I82: 4622.143 -> 4725.8335
$ws.Cells.Item(82,11).Value = 4725.8335  # K82: 4622.143 -> 4725.8335
$ws.Cells.Item(82,13).Value = -4342.8335  # M82: -4239.143 -> -4342.8335

$ws.Cells.Item(85,8).Value = 14336.429  # H85: 13044.375 -> 14336.429
$ws.Cells.Item(85,9).Value = 4725.8335  # I85: 4622.143 -> 4725.8335
$ws.Cells.Item(85,11).Value = 4725.8335  # K85: 4622.143 -> 4725.8335
$ws.Cells.Item(85,13).Value = -3399.8335  # M85: -3296.143 -> -3399.8335

$ws.Cells.Item(107,8).Value = 14000  # H107: 3022.4285 -> 14000
$ws.Cells.Item(107,9).Value = 0  # I107: 1192.8334 -> 0
$ws.Cells.Item(107,11).Value = 0  # K107: 1192.8334 -> 0
$ws.Cells.Item(107,13).ClearContents()  # M107: 727.1666 -> (removed)

$ws.Cells.Item(134,8).Value = 2453.973  # H134: 2501.5405 -> 2453.973
$ws.Cells.Item(134,9).Value = 2110.1667  # I134: 2140.1724 -> 2110.1667
$ws.Cells.Item(134,10).Value = 3927.4285  # J134: 3811.5 -> 3927.4285
$ws.Cells.Item(134,11).Value = 6330.500100000001  # K134: 6420.5172 -> 6330.500100000001
$ws.Cells.Item(134,12).Value = 11782.2855  # L134: 11434.5 -> 11782.2855
$ws.Cells.Item(134,13).Value = -3795.500100000001  # M134: -3885.5172 -> -3795.500100000001
$ws.Cells.Item(134,14).Value = -16852.2855  # N134: -16504.5 -> -16852.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22,8).Value = 462.375  # H22: 476.53333 -> 462.375
$ws.Cells.Item(22,10).Value = 549.75  # J22: 592.5714 -> 549.75
$ws.Cells.Item(22,12).Value = 549.75  # L22: 592.5714 -> 549.75
$ws.Cells.Item(22,14).Value = -1249.75  # N22: -1292.5714 -> -1249.75

$ws.Cells.Item(58,8).Value = 2242.818  # H58: 2316 -> 2242.818
$ws.Cells.Item(58,9).Value = 2002.75  # I58: 2073 -> 2002.75
$ws.Cells.Item(58,11).Value = 2002.75  # K58: 2073 -> 2002.75
$ws.Cells.Item(58,13).Value = -1799.75  # M58: -1870 -> -1799.75

$ws.Cells.Item(59,8).Value = 0  # H59: 32499.5 -> 0
$ws.Cells.Item(59,10).Value = 0  # J59: 32499.5 -> 0
$ws.Cells.Item(59,12).Value = 0  # L59: 32499.5 -> 0
$ws.Cells.Item(59,14).ClearContents()  # N59: -34789.5 -> (removed)

$ws.Cells.Item(64,8).Value = 48649.152  # H64: 0 -> 48649.152
$ws.Cells.Item(64,10).Value = 48649.152  # J64: 0 -> 48649.152
$ws.Cells.Item(64,12).Value = 48649.152  # L64: 0 -> 48649.152
$ws.Cells.Item(64,14).Value = -49145.152  # N64: None -> -49145.152

$ws.Cells.Item(67,8).Value = 48649.152  # H67: 0 -> 48649.152
$ws.Cells.Item(67,10).Value = 48649.152  # J67: 0 -> 48649.152
$ws.Cells.Item(67,12).Value = 48649.152  # L67: 0 -> 48649.152
$ws.Cells.Item(67,14).Value = -50365.152  # N67: None -> -50365.152

$ws.Cells.Item(68,8).Value = 59332.668  # H68: 62749.5 -> 59332.668
$ws.Cells.Item(68,10).Value = 79998  # J68: 76499 -> 79998
$ws.Cells.Item(68,12).Value = 79998  # L68: 76499 -> 79998
$ws.Cells.Item(68,14).Value = -81496  # N68: -77997 -> -81496

$ws.Cells.Item(71,8).Value = 59332.668  # H71: 62749.5 -> 59332.668
$ws.Cells.Item(71,10).Value = 79998  # J71: 76499 -> 79998
$ws.Cells.Item(71,12).Value = 239994  # L71: 229497 -> 239994
$ws.Cells.Item(71,14).Value = -247482  # N71: -236985 -> -247482

$ws.Cells.Item(107,8).Value = 1718.8462  # H107: 1164.7 -> 1718.8462
$ws.Cells.Item(107,9).Value = 1585.5  # I107: 1046.7059 -> 1585.5
$ws.Cells.Item(107,10).Value = 2163.3333  # J107: 1833.3334 -> 2163.3333
$ws.Cells.Item(107,11).Value = 1585.5  # K107: 1046.7059 -> 1585.5
$ws.Cells.Item(107,12).Value = 2163.3333  # L107: 1833.3334 -> 2163.3333
$ws.Cells.Item(107,13).Value = 334.5  # M107: 873.2941000000001 -> 334.5
$ws.Cells.Item(107,14).Value = -6003.3333  # N107: -5673.3334 -> -6003.3333

$ws.Cells.Item(132,8).Value = 2015.9615  # H132: 1971.5 -> 2015.9615
$ws.Cells.Item(132,9).Value = 1769.1818  # I132: 1801.0952 -> 1769.1818
$ws.Cells.Item(132,10).Value = 3373.25  # J132: 3164.3333 -> 3373.25
$ws.Cells.Item(132,11).Value = 5307.5454  # K132: 5403.2856 -> 5307.5454
$ws.Cells.Item(132,12).Value = 10119.75  # L132: 9492.999899999999 -> 10119.75
$ws.Cells.Item(132,13).Value = -2777.5454  # M132: -2873.2856 -> -2777.5454
$ws.Cells.Item(132,14).Value = -15179.75  # N132: -14552.9999 -> -15179.75

$ws.Cells.Item(134,8).Value = 2065.375  # H134: 2238.9524 -> 2065.375
$ws.Cells.Item(134,9).Value = 2081.818  # I134: 2276.2632 -> 2081.818
$ws.Cells.Item(134,11).Value = 6245.454000000001  # K134: 6828.7896 -> 6245.454000000001
$ws.Cells.Item(134,13).Value = -3710.454000000001  # M134: -4293.7896 -> -3710.454000000001

$ws.Cells.Item(136,8).Value = 2242.818  # H136: 2316 -> 2242.818
$ws.Cells.Item(136,9).Value = 2002.75  # I136: 2073 -> 2002.75
$ws.Cells.Item(136,11).Value = 6008.25  # K136: 6219 -> 6008.25
$ws.Cells.Item(136,13).Value = -3458.25  # M136: -3669 -> -3458.25

$ws.Cells.Item(141,8).Value = 428531  # H141: 412433.62 -> 428531
$ws.Cells.Item(141,9).Value = 74999.664  # I141: 58749.5 -> 74999.664
$ws.Cells.Item(141,11).Value = 74999.664  # K141: 58749.5 -> 74999.664
$ws.Cells.Item(141,13).Value = -69819.664  # M141: -53569.5 -> -69819.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2,8).Value = 1384.125  # H2: 1149.8 -> 1384.125
$ws.Cells.Item(2,9).Value = 105.166664  # I2: 128.875 -> 105.166664
$ws.Cells.Item(2,10).Value = 2151.5  # J2: 1830.4166 -> 2151.5
$ws.Cells.Item(2,11).Value = 630.999984  # K2: 773.25 -> 630.999984
$ws.Cells.Item(2,12).Value = 12909  # L2: 10982.4996 -> 12909
$ws.Cells.Item(2,13).Value = -517.999984  # M2: -660.25 -> -517.999984
$ws.Cells.Item(2,14).Value = -13135  # N2: -11208.4996 -> -13135

$ws.Cells.Item(5,8).Value = 1302.75  # H5: 1341.2727 -> 1302.75
$ws.Cells.Item(5,9).Value = 726.4286  # I5: 729.2857 -> 726.4286
$ws.Cells.Item(5,10).Value = 2109.6  # J5: 2412.25 -> 2109.6
$ws.Cells.Item(5,11).Value = 2179.2858  # K5: 2187.8571 -> 2179.2858
$ws.Cells.Item(5,12).Value = 6328.799999999999  # L5: 7236.75 -> 6328.799999999999
$ws.Cells.Item(5,13).Value = -2067.2858  # M5: -2075.8571 -> -2067.2858
$ws.Cells.Item(5,14).Value = -6552.799999999999  # N5: -7460.75 -> -6552.799999999999

$ws.Cells.Item(34,8).Value = 2778.8333  # H34: 1025.375 -> 2778.8333
$ws.Cells.Item(34,9).Value = 224.33333  # I34: 171.85715 -> 224.33333
$ws.Cells.Item(34,10).Value = 5333.3335  # J34: 7000 -> 5333.3335
$ws.Cells.Item(34,11).Value = 672.99999  # K34: 515.5714499999999 -> 672.99999
$ws.Cells.Item(34,12).Value = 16000.0005  # L34: 21000 -> 16000.0005
$ws.Cells.Item(34,13).Value = -588.99999  # M34: -431.5714499999999 -> -588.99999
$ws.Cells.Item(34,14).Value = -16168.0005  # N34: -21168 -> -16168.0005

$ws.Cells.Item(113,8).Value = 3217.1738  # H113: 3086.8262 -> 3217.1738
$ws.Cells.Item(113,9).Value = 0  # I113: 1999.5 -> 0
$ws.Cells.Item(113,10).Value = 3217.1738  # J113: 3190.3809 -> 3217.1738
$ws.Cells.Item(113,11).Value = 0  # K113: 5998.5 -> 0
$ws.Cells.Item(113,12).Value = 9651.5214  # L113: 9571.1427 -> 9651.5214
$ws.Cells.Item(113,13).ClearContents()  # M113: -3828.5 -> (removed)
$ws.Cells.Item(113,14).Value = -13991.5214  # N113: -13911.1427 -> -13991.5214

$ws.Cells.Item(118,8).Value = 9579.412  # H118: 9326.611000000001 -> 9579.412
$ws.Cells.Item(118,9).Value = 8952.375  # I118: 8516.444 -> 8952.375
$ws.Cells.Item(118,11).Value = 26857.125  # K118: 25549.332 -> 26857.125
$ws.Cells.Item(118,13).Value = -25614.125  # M118: -24306.332 -> -25614.125

$ws.Cells.Item(132,8).Value = 3127.389  # H132: 2916.524 -> 3127.389
$ws.Cells.Item(132,9).Value = 1645.8182  # I132: 1631.3846 -> 1645.8182
$ws.Cells.Item(132,10).Value = 5455.5713  # J132: 5004.875 -> 5455.5713
$ws.Cells.Item(132,11).Value = 14812.3638  # K132: 14682.4614 -> 14812.3638
$ws.Cells.Item(132,12).Value = 49100.14169999999  # L132: 45043.875 -> 49100.14169999999
$ws.Cells.Item(132,13).Value = -12282.3638  # M132: -12152.4614 -> -12282.3638
$ws.Cells.Item(132,14).Value = -54160.14169999999  # N132: -50103.875 -> -54160.14169999999

$ws.Cells.Item(135,8).Value = 1302.75  # H135: 1341.2727 -> 1302.75
$ws.Cells.Item(135,9).Value = 726.4286  # I135: 729.2857 -> 726.4286
$ws.Cells.Item(135,10).Value = 2109.6  # J135: 2412.25 -> 2109.6
$ws.Cells.Item(135,11).Value = 6537.8574  # K135: 6563.571300000001 -> 6537.8574
$ws.Cells.Item(135,12).Value = 18986.4  # L135: 21710.25 -> 18986.4
$ws.Cells.Item(135,13).Value = -4002.8574  # M135: -4028.571300000001 -> -4002.8574
$ws.Cells.Item(135,14).Value = -24056.4  # N135: -26780.25 -> -24056.4

$ws.Cells.Item(136,8).Value = 5902.0303  # H136: 5163.028 -> 5902.0303
$ws.Cells.Item(136,9).Value = 1791.0714  # I136: 1703 -> 1791.0714
$ws.Cells.Item(136,10).Value = 8931.157999999999  # J136: 9030.117 -> 8931.157999999999
$ws.Cells.Item(136,11).Value = 5373.2142  # K136: 5109 -> 5373.2142
$ws.Cells.Item(136,12).Value = 26793.474  # L136: 27090.351 -> 26793.474
$ws.Cells.Item(136,13).Value = -273.2142000000003  # M136: -9 -> -273.2142000000003
$ws.Cells.Item(136,14).Value = -36993.474  # N136: -37290.351 -> -36993.474

$ws.Cells.Item(137,8).Value = 53575492  # H137: 57696430 -> 53575492
$ws.Cells.Item(137,10).Value = 4741.25  # J137: 4872.364 -> 4741.25
$ws.Cells.Item(137,12).Value = 14223.75  # L137: 14617.092 -> 14223.75
$ws.Cells.Item(137,14).Value = -24423.75  # N137: -24817.092 -> -24423.75

$ws.Cells.Item(138,8).Value = 4374.625  # H138: 4374.75 -> 4374.625
$ws.Cells.Item(138,9).Value = 4374.625  # I138: 4374.75 -> 4374.625
$ws.Cells.Item(138,11).Value = 13123.875  # K138: 13124.25 -> 13123.875
$ws.Cells.Item(138,13).Value = -7983.875  # M138: -7984.25 -> -7983.875

$ws.Cells.Item(140,8).Value = 5392.593  # H140: 5195.8276 -> 5392.593
$ws.Cells.Item(140,9).Value = 3180.55  # I140: 3122.2727 -> 3180.55
$ws.Cells.Item(140,11).Value = 9541.650000000001  # K140: 9366.8181 -> 9541.650000000001
$ws.Cells.Item(140,13).Value = -4361.650000000001  # M140: -4186.8181 -> -4361.650000000001

$ws.Cells.Item(141,8).Value = 9568.058999999999  # H141: 9573.235000000001 -> 9568.058999999999
$ws.Cells.Item(141,9).Value = 8952.714  # I141: 8965.286 -> 8952.714
$ws.Cells.Item(141,11).Value = 26858.142  # K141: 26895.858 -> 26858.142
$ws.Cells.Item(141,13).Value = -21678.142  # M141: -21715.858 -> -21678.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11,8).Value = 4696077  # H11: 6669750 -> 4696077
$ws.Cells.Item(11,9).Value = 2337889.2  # I11: 2503875 -> 2337889.2
$ws.Cells.Item(11,10).Value = 10002000  # J11: 15001500 -> 10002000
$ws.Cells.Item(11,11).Value = 2337889.2  # K11: 2503875 -> 2337889.2
$ws.Cells.Item(11,12).Value = 10002000  # L11: 15001500 -> 10002000
$ws.Cells.Item(11,13).Value = -2337750.2  # M11: -2503736 -> -2337750.2
$ws.Cells.Item(11,14).Value = -10002278  # N11: -15001778 -> -10002278

$ws.Cells.Item(43,8).Value = 45650  # H43: 45574.5 -> 45650
$ws.Cells.Item(43,10).Value = 45650  # J43: 45574.5 -> 45650
$ws.Cells.Item(43,12).Value = 45650  # L43: 45574.5 -> 45650
$ws.Cells.Item(43,14).Value = -45952  # N43: -45876.5 -> -45952

$ws.Cells.Item(46,8).Value = 35067.5  # H46: 34588.047 -> 35067.5
$ws.Cells.Item(46,10).Value = 57083.332  # J46: 54615.31 -> 57083.332
$ws.Cells.Item(46,12).Value = 57083.332  # L46: 54615.31 -> 57083.332
$ws.Cells.Item(46,14).Value = -57395.332  # N46: -54927.31 -> -57395.332

$ws.Cells.Item(57,8).Value = 21999  # H57: 21998.75 -> 21999
$ws.Cells.Item(57,10).Value = 21999  # J57: 21998.75 -> 21999
$ws.Cells.Item(57,12).Value = 21999  # L57: 21998.75 -> 21999
$ws.Cells.Item(57,14).Value = -23639  # N57: -23638.75 -> -23639

$ws.Cells.Item(107,8).Value = 377.4  # H107: 310.9 -> 377.4
$ws.Cells.Item(107,9).Value = 377.4  # I107: 298 -> 377.4
$ws.Cells.Item(107,10).Value = 0  # J107: 341 -> 0
$ws.Cells.Item(107,11).Value = 377.4  # K107: 298 -> 377.4
$ws.Cells.Item(107,12).Value = 0  # L107: 341 -> 0
$ws.Cells.Item(107,13).Value = 1542.6  # M107: 1622 -> 1542.6
$ws.Cells.Item(107,14).ClearContents()  # N107: -4181 -> (removed)

$ws.Cells.Item(113,8).Value = 1998  # H113: 1999 -> 1998
$ws.Cells.Item(113,10).Value = 1997.5  # J113: 2000 -> 1997.5
$ws.Cells.Item(113,12).Value = 1997.5  # L113: 2000 -> 1997.5
$ws.Cells.Item(113,14).Value = -6337.5  # N113: -6340 -> -6337.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20,8).Value = 1622369  # H20: 1622383.4 -> 1622369
$ws.Cells.Item(20,10).Value = 6176581.5  # J20: 6176636.5 -> 6176581.5
$ws.Cells.Item(20,12).Value = 6176581.5  # L20: 6176636.5 -> 6176581.5
$ws.Cells.Item(20,14).Value = -6177033.5  # N20: -6177088.5 -> -6177033.5

$ws.Cells.Item(132,8).Value = 6215.3706  # H132: 5990.3486 -> 6215.3706
$ws.Cells.Item(132,9).Value = 5582.8906  # I132: 5234.541 -> 5582.8906
$ws.Cells.Item(132,11).Value = 16748.6718  # K132: 15703.623 -> 16748.6718
$ws.Cells.Item(132,13).Value = -14218.6718  # M132: -13173.623 -> -14218.6718

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54,8).Value = 39500  # H54: 41332.668 -> 39500
$ws.Cells.Item(54,10).Value = 39000  # J54: 41999 -> 39000
$ws.Cells.Item(54,12).Value = 39000  # L54: 41999 -> 39000
$ws.Cells.Item(54,14).Value = -40040  # N54: -43039 -> -40040

$ws.Cells.Item(122,8).Value = 5450.222  # H122: 7103.643 -> 5450.222
$ws.Cells.Item(122,9).Value = 5033.2  # I122: 6879.4287 -> 5033.2
$ws.Cells.Item(122,10).Value = 5971.5  # J122: 7327.857 -> 5971.5
$ws.Cells.Item(122,11).Value = 15099.6  # K122: 20638.2861 -> 15099.6
$ws.Cells.Item(122,12).Value = 17914.5  # L122: 21983.571 -> 17914.5
$ws.Cells.Item(122,13).Value = -12649.6  # M122: -18188.2861 -> -12649.6
$ws.Cells.Item(122,14).Value = -22814.5  # N122: -26883.571 -> -22814.5

$ws.Cells.Item(132,8).Value = 6670115.5  # H132: 7095812 -> 6670115.5
$ws.Cells.Item(132,9).Value = 14495544  # I132: 16669746 -> 14495544
$ws.Cells.Item(132,11).Value = 43486632  # K132: 50009238 -> 43486632
$ws.Cells.Item(132,13).Value = -43484102  # M132: -50006708 -> -43484102

$ws.Cells.Item(136,8).Value = 5770.1665  # H136: 5858.8936 -> 5770.1665
$ws.Cells.Item(136,9).Value = 4044.9395  # I136: 4121.3438 -> 4044.9395
$ws.Cells.Item(136,11).Value = 12134.8185  # K136: 12364.0314 -> 12134.8185
$ws.Cells.Item(136,13).Value = -9584.818499999999  # M136: -9814.0314 -> -9584.818499999999
